# cs-en-us-075pct.xlsx -- 75th Precinct weekly CompStat report refresh
# New crime data collected: advance the report to Volume 30 Number 38,
# covering the week 9/18/2023 through 9/24/2023, and update every
# Week to Date / 28 Day / Year to Date / 2 Year complaint figure below.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Masthead: report volume/number and the week-covering date range
$ws.Range("A8").Value = "Volume 30   Number  38"
$ws.Range("C9").Value = "Report Covering the Week  9/18/2023  Through  9/24/2023"

# Row 14 - Murder
$ws.Range("F14").Value = 1
$ws.Range("H14").Value = -80
$ws.Range("J14").Value = 17
$ws.Range("K14").Value = 0
$ws.Range("N14").Value = -81.720430107526

# Row 15 - Rape
$ws.Range("C15").Value = 3
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 50
$ws.Range("F15").Value = 8
$ws.Range("H15").Value = 14.285714285714
$ws.Range("I15").Value = 49
$ws.Range("J15").Value = 39
$ws.Range("K15").Value = 25.641025641025
$ws.Range("L15").Value = 16.666666666666
$ws.Range("M15").Value = 8.888888888888
$ws.Range("N15").Value = -45.555555555555

# Row 16 - Robbery
$ws.Range("C16").Value = 16
$ws.Range("E16").Value = 23.076923076923
$ws.Range("F16").Value = 58
$ws.Range("G16").Value = 59
$ws.Range("H16").Value = -1.694915254237
$ws.Range("I16").Value = 500
$ws.Range("J16").Value = 570
$ws.Range("K16").Value = -12.280701754386
$ws.Range("L16").Value = 15.207373271889
$ws.Range("M16").Value = -6.716417910447
$ws.Range("N16").Value = -76.841130152848

# Row 17 - Fel. Assault
$ws.Range("D17").Value = 18
$ws.Range("E17").Value = 27.777777777777
$ws.Range("F17").Value = 78
$ws.Range("G17").Value = 89
$ws.Range("H17").Value = -12.359550561797
$ws.Range("I17").Value = 787
$ws.Range("J17").Value = 823
$ws.Range("K17").Value = -4.374240583232
$ws.Range("L17").Value = 20.705521472392
$ws.Range("M17").Value = 39.045936395759
$ws.Range("N17").Value = -30.538393645189

# Row 18 - Burglary
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = 16.666666666666
$ws.Range("F18").Value = 20
$ws.Range("G18").Value = 27
$ws.Range("H18").Value = -25.925925925925
$ws.Range("I18").Value = 271
$ws.Range("J18").Value = 294
$ws.Range("K18").Value = -7.823129251700
$ws.Range("L18").Value = -5.574912891986
$ws.Range("M18").Value = -19.822485207100
$ws.Range("N18").Value = -80.014749262536

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 27
$ws.Range("D19").Value = 30
$ws.Range("E19").Value = -10
$ws.Range("F19").Value = 86
$ws.Range("G19").Value = 111
$ws.Range("H19").Value = -22.522522522522
$ws.Range("I19").Value = 800
$ws.Range("J19").Value = 988
$ws.Range("K19").Value = -19.028340080971
$ws.Range("L19").Value = 7.671601615074
$ws.Range("M19").Value = 57.170923379174
$ws.Range("N19").Value = 12.517580872011

# Row 20 - G.L.A.
$ws.Range("C20").Value = 13
$ws.Range("D20").Value = 11
$ws.Range("E20").Value = 18.181818181818
$ws.Range("G20").Value = 40
$ws.Range("H20").Value = 10
$ws.Range("I20").Value = 399
$ws.Range("J20").Value = 355
$ws.Range("K20").Value = 12.394366197183
$ws.Range("L20").Value = 3.636363636363
$ws.Range("M20").Value = 76.548672566371
$ws.Range("N20").Value = -79.283489096573

# Row 21 - TOTAL
$ws.Range("C21").Value = 89
$ws.Range("D21").Value = 81
$ws.Range("E21").Value = 9.876543209876
$ws.Range("F21").Value = 295
$ws.Range("G21").Value = 338
$ws.Range("H21").Value = -12.721893491124
$ws.Range("I21").Value = 2823
$ws.Range("J21").Value = 3086
$ws.Range("K21").Value = -8.522359040829
$ws.Range("L21").Value = 10.230378758297
$ws.Range("M21").Value = 25.802139037433
$ws.Range("N21").Value = -62.198714515265

# Row 22 - Transit
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = -16.666666666666
$ws.Range("I22").Value = 34
$ws.Range("J22").Value = 59
$ws.Range("K22").Value = -42.372881355932
$ws.Range("L22").Value = 13.333333333333
$ws.Range("M22").Value = -27.659574468085

# Row 23 - Housing
$ws.Range("C23").Value = 6
$ws.Range("E23").Value = -25
$ws.Range("F23").Value = 20
$ws.Range("G23").Value = 32
$ws.Range("H23").Value = -37.5
$ws.Range("I23").Value = 245
$ws.Range("J23").Value = 265
$ws.Range("K23").Value = -7.547169811320
$ws.Range("L23").Value = 15.566037735849
$ws.Range("M23").Value = 47.590361445783

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 59
$ws.Range("D24").Value = 62
$ws.Range("E24").Value = -4.838709677419
$ws.Range("F24").Value = 182
$ws.Range("G24").Value = 228
$ws.Range("H24").Value = -20.175438596491
$ws.Range("I24").Value = 1743
$ws.Range("J24").Value = 2140
$ws.Range("K24").Value = -18.551401869158
$ws.Range("L24").Value = 7.526218383713
$ws.Range("M24").Value = 60.941828254847

# Row 25 - Misd. Assault
$ws.Range("D25").Value = 16
$ws.Range("E25").Value = 50
$ws.Range("F25").Value = 95
$ws.Range("G25").Value = 91
$ws.Range("H25").Value = 4.395604395604
$ws.Range("I25").Value = 923
$ws.Range("J25").Value = 950
$ws.Range("K25").Value = -2.842105263157
$ws.Range("L25").Value = 26.438356164383
$ws.Range("M25").Value = -34.677990092002

# Row 26 - UCR Rape*
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = 50
$ws.Range("F26").Value = 10
$ws.Range("H26").Value = 25
$ws.Range("I26").Value = 73
$ws.Range("J26").Value = 57
$ws.Range("K26").Value = 28.070175438596
$ws.Range("L26").Value = 0

# Row 27 - Other Sex Crimes
$ws.Range("C27").Value = 5
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = 150
$ws.Range("F27").Value = 10
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = 42.857142857142
$ws.Range("I27").Value = 90
$ws.Range("J27").Value = 91
$ws.Range("K27").Value = -1.098901098901
$ws.Range("L27").Value = 0

# Row 28 - Shooting Vic.
$ws.Range("D28").Value = 2
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 9
$ws.Range("H28").Value = -88.888888888888
$ws.Range("J28").Value = 73
$ws.Range("K28").Value = -34.246575342465
$ws.Range("M28").Value = -43.529411764705
$ws.Range("N28").Value = -81.954887218045

# Row 29 - Shooting Inc.
$ws.Range("D29").Value = 2
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 8
$ws.Range("H29").Value = -87.5
$ws.Range("J29").Value = 64
$ws.Range("K29").Value = -34.375
$ws.Range("M29").Value = -40
$ws.Range("N29").Value = -82.426778242677

# Row 30 - Hate Crimes
$ws.Range("L30").Value = -60
